$wb = $excel.ActiveWorkbook

# --- Sheet2: give column A an explicit width (~40.16 chars) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns.Item(1).ColumnWidth = 39.33

# --- Sheet3: filenames table updates ---
$ws3 = $wb.Worksheets.Item("Sheet3")

# Update the filename referenced by A2 (new CoolTerm capture log)
$ws3.Range("A2").Value = "../data/CoolTerm Capture 2023-02-02 10-56-19.txt"

# Remove bold from the header row (A1:F1)
$ws3.Range("A1:F1").Font.Bold = $false

# Give column A an explicit width (41 chars)
$ws3.Columns.Item(1).ColumnWidth = 40.17

# Update numeric data B2:E5 with the new training-model sample values
$ws3.Range("B2").Value = 356
$ws3.Range("C2").Value = 1314
$ws3.Range("D2").Value = 2669
$ws3.Range("E2").Value = 2008

$ws3.Range("B3").Value = 620
$ws3.Range("C3").Value = 1512
$ws3.Range("D3").Value = 2769
$ws3.Range("E3").Value = 2273

$ws3.Range("B4").Value = 818
$ws3.Range("C4").Value = 1678
$ws3.Range("D4").Value = 2967
$ws3.Range("E4").Value = 2372

$ws3.Range("B5").Value = 1182
$ws3.Range("C5").Value = 1843
$ws3.Range("D5").Value = 3033
$ws3.Range("E5").Value = 2471

# Move the selection on Sheet3 to D6:E7
$ws3.Activate()
$ws3.Range("D6:E7").Select()
